$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the last row of the time log (week 16) with hours and activity.
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Meetings, fixes and poster"

# Move the active selection to F22, as in the edited workbook.
$ws.Range("F22").Select()
